$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo in the marketing amount header and select new active cell
$ws.Range("C1").Value = "Marketing_Amount"
$ws.Range("D1").Value = "Num_Email_Subscribers"

$ws.Range("C4").Select()
